# Edit the map workbook:
#  - rename "seller_room" cells to "trader_room"
#  - rename "tall_man_room" cells to "long_hallway"
#  - rename "little_boy_room" cell to "little_boy_bossroom"
#  - move the active selection from K6 to C3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value2
    if ($val -eq "seller_room") {
        $cell.Value = "trader_room"
    } elseif ($val -eq "tall_man_room") {
        $cell.Value = "long_hallway"
    } elseif ($val -eq "little_boy_room") {
        $cell.Value = "little_boy_bossroom"
    }
}

$ws.Range("C3").Select() | Out-Null
